$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-4: "Did Harvest Occur?" changes from Yes to No,
# Species (F) is cleared, and Unknown Sex Count (J) changes from 1 to 0
for ($row = 2; $row -le 4; $row++) {
    $ws.Cells.Item($row, 2).Value = "No"
    $ws.Cells.Item($row, 6).Value = ""
    $ws.Cells.Item($row, 10).Value = 0
}
